$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44, shifting existing rows 44:103 down to 45:104
$ws.Rows.Item(44).Insert()

# Populate the new row 44 with this week's data
$ws.Range("A44").Value = 4
$ws.Range("B44").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C44").Value = "Los Lagos"
$ws.Range("D44").Value = 44467
$ws.Range("D44").NumberFormat = $ws.Range("D45").NumberFormat
$ws.Range("E44").Value = 10
$ws.Range("F44").Value = "Fruta"
$ws.Range("G44").Value = 100102
$ws.Range("H44").Value = "Cítricos"
$ws.Range("I44").Value = 100102004
$ws.Range("J44").Value = "Mandarina"
$ws.Range("K44").Value = "Murcott"
$ws.Range("L44").Value = "Primera"
$ws.Range("M44").Value = 600
$ws.Range("N44").Value = 6500
$ws.Range("O44").Value = 6500
$ws.Range("P44").Value = 6500
$ws.Range("Q44").Value = "$/bandeja 10 kilos"
$ws.Range("R44").Value = "Provincia de Limarí"
$ws.Range("S44").Value = 650
$ws.Range("T44").Value = 10
